$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row before row 48. This pushes the existing rows 48-54
#    (empty spacer, "Quelle: Mikrozensus", blank spacer, ministry credit,
#    copyright line, reproduction notice, integrationsmonitoring URL) down
#    to rows 49-55.
# ---------------------------------------------------------------------------
$ws.Rows.Item(48).Insert()

# ---------------------------------------------------------------------------
# 2) Footnote 4) in B47 referenced a long URL inline. Split it: B47 keeps the
#    lead-in text (through "Bundesamtes: ") and the URL itself moves into the
#    freshly inserted row 48 as its own hyperlinked cell.
# ---------------------------------------------------------------------------
$marker = "Bundesamtes: "
$full = $ws.Range("B47").Value2
$cut = $full.IndexOf($marker) + $marker.Length
$headText = $full.Substring(0, $cut)
$urlText = $full.Substring($cut)

$ws.Range("B47").Value2 = $headText

$ws.Range("B48").Value2 = $urlText
$ws.Hyperlinks.Add($ws.Range("B48"), $urlText) | Out-Null

# ---------------------------------------------------------------------------
# 3) Row 9 (the "Insgesamt" / eigener & ohne eigener Migrationserfahrung sub
#    header row) now wraps onto multiple lines, so it needs extra height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 4) Widen the three "eigener/ohne eigene Migrationserfahrung" column pairs
#    (G:H, L:M, Q:R) so the new wrapped header text in row 9 fits better.
# ---------------------------------------------------------------------------
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 13.86
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 13.86
$ws.Range("Q1:R1").EntireColumn.ColumnWidth = 13.86

$wb.Save()
